$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2
$wdFindContinue = 1
$wdReplaceAll = 2

# 1. Title paragraph: "1 What is Blender " -> "7 Relax, Ummm... Err... Smooth Vertices"
$pTitle = $d.Paragraphs(1).Range
$pTitle.Find.Execute("1 What is Blender ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "7 Relax, Ummm... Err... Smooth Vertices", $wdReplaceAll)

# 2. Keywords paragraph: insert "Relax, Smooth, " before "Blender, 3D Modeling, Animation, Graphic Art"
$pKeywords = $d.Paragraphs(6).Range
$pKeywords.Find.Execute(" Blender, 3D Modeling, Animation, Graphic Art ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, " Relax, Smooth, Blender, 3D Modeling, Animation, Graphic Art ", $wdReplaceAll)

# 3. Description paragraph: replace tail of description
$pDescription = $d.Paragraphs(9).Range
$pDescription.Find.Execute("what the 3D modeling program ""Blender "" is all about./", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "two different ways to relax or smooth vertices in Blender/", $wdReplaceAll)

# 4. Category paragraph: insert "Relax, Smooth, " before "Blender, 3D Modeling, Animation, Graphic Art"
$pCategory = $d.Paragraphs(11).Range
$pCategory.Find.Execute("Blender, 3D Modeling, Animation, Graphic Art", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Relax, Smooth, Blender, 3D Modeling, Animation, Graphic Art", $wdReplaceAll)

# 5. Revised paragraph: update date
$pRevised = $d.Paragraphs(18).Range
$pRevised.Find.Execute("Wednesday, December 11, 2024", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Sunday, February 23, 2025", $wdReplaceAll)

# 6. Url paragraph: update url
$pUrl = $d.Paragraphs(20).Range
$pUrl.Find.Execute("Enlightenment/Articles/2024/8-Blender-2024/1-What-Is-Blender/1-What-Is-Blender.html", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "Enlightenment/Articles/2025/1-Blender-Continued/7-Loop-Tools/7-Relax-Umm-Err-Smooth-Vertices/7-Relax-Ummm-Err-Smooth-Vertices.html", $wdReplaceAll)
